$wb = $excel.ActiveWorkbook

# --- Sheet "ChosenTreatment" ---
$ws1 = $wb.Worksheets.Item("ChosenTreatment")
$ws1.Range("A2").Value = "{'SexualOrientation': 'Straight or heterosexual'}"
$ws1.Range("B2").Value = "{'DevType': 'Back-end developer'}"

# --- Sheet "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 26

# --- Sheet "Subgroups" ---
$ws3 = $wb.Worksheets.Item("Subgroups")

$data = @(
  @("{'Gender': '1'}", 33253, 5920.984439082945, -792.3712621303766),
  @("{'Student': '1'}", 29526, 5617.758312088003, -1095.597389125319),
  @("{'Hobby': '1'}", 28842, 5516.21828962618, -1197.137411587141),
  @("{'HDI': '1'}", 27581, 8280.034330059932, 1566.678628846611),
  @("{'RaceEthnicity': '1'}", 27379, 6509.936391650183, -203.4193095631381),
  @("{'Dependents': '2'}", 24167, 7072.162954162452, 358.8072529491301),
  @("{'UndergradMajor': '2'}", 22173, 3542.149042873403, -3171.206658339918),
  @("{'Gender': '1', 'Student': '1'}", 27820, 4860.761884746209, -1852.593816467112),
  @("{'Gender': '1', 'Hobby': '1'}", 27456, 5084.527001737077, -1628.828699476245),
  @("{'Hobby': '1', 'Student': '1'}", 23894, 4208.446115890501, -2504.90958532282),
  @("{'Gender': '1', 'Hobby': '1', 'Student': '1'}", 22758, 4081.320153739394, -2632.035547473928),
  @("{'Gender': '1', 'HDI': '1'}", 25944, 7219.315451383841, 505.9597501705193),
  @("{'HDI': '1', 'Student': '1'}", 23773, 7095.77558069452, 382.4198794811982),
  @("{'Hobby': '1', 'HDI': '1'}", 22660, 7025.684796894849, 312.3290956815272),
  @("{'Gender': '1', 'HDI': '1', 'Student': '1'}", 22370, 6176.689060932703, -536.6666402806186),
  @("{'Gender': '1', 'Hobby': '1', 'HDI': '1'}", 21574, 6412.706957617019, -300.6487435963027),
  @("{'Gender': '1', 'RaceEthnicity': '1'}", 25910, 6471.012774252379, -242.3429269609423),
  @("{'RaceEthnicity': '1', 'Student': '1'}", 23264, 5618.412897795355, -1094.942803417966),
  @("{'RaceEthnicity': '1', 'HDI': '1'}", 23146, 7669.000695259265, 955.6449940459433),
  @("{'RaceEthnicity': '1', 'Hobby': '1'}", 22529, 5345.204113474854, -1368.151587738467),
  @("{'Gender': '1', 'RaceEthnicity': '1', 'Student': '1'}", 22014, 5982.017555487582, -731.3381457257392),
  @("{'Gender': '1', 'RaceEthnicity': '1', 'HDI': '1'}", 21900, 7566.667087338407, 853.3113861250858),
  @("{'RaceEthnicity': '1', 'HDI': '1', 'Student': '1'}", 20128, 6694.06933007496, -19.28637113836157),
  @("{'Gender': '1', 'RaceEthnicity': '1', 'Hobby': '1'}", 21544, 5551.102665918384, -1162.253035294938),
  @("{'Gender': '1', 'Dependents': '2'}", 22568, 6660.009269424866, -53.34643178845545),
  @("{'Gender': '1', 'UndergradMajor': '2'}", 21070, 2140.616974317039, -4572.738726896283)
)

$row = 2
foreach ($r in $data) {
    $ws3.Cells.Item($row, 1).Value = $r[0]
    $ws3.Cells.Item($row, 2).Value = $r[1]
    $ws3.Cells.Item($row, 3).Value = $r[2]
    $ws3.Cells.Item($row, 4).Value = $r[3]
    $row++
}
